# Update the business_type (column B) classification for a handful of
# industries so they map to the standard category labels instead of the
# stray lower-case duplicates ("automotive", "pet related") or the
# mismatched category that used to sit there.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = "professional services"   # Automotive
$ws.Range("B3").Value  = "Retail"                  # Food & Restaurants
$ws.Range("B6").Value  = "professional services"   # Building & Storage
$ws.Range("B7").Value  = "Retail"                  # Pet Related
$ws.Range("B38").Value = "Retail"                  # Frozen Yogurt / Ice Cream
$ws.Range("B39").Value = "Retail"                  # Car Rental
$ws.Range("B40").Value = "Retail"                  # Car Wash

# Selection moved to E11 when the file was last saved.
$ws.Range("E11").Select()
